$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before row 11 (the old "Budget View" row), shifting rows 11-21 down to 12-22.
$ws.Rows("11:11").Insert()

# New row 11 should read "Charge Type" in column A (column B left blank),
# matching the style of the surrounding "first group" label rows (3-12).
$ws.Range("A11").Value = "Charge Type"
